$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates in D2:D3 (currently 2022-01-13 / serial 44574) with D4:D5
# (currently 2021-12-29 / serial 44559), per the weekly Fruta/Hortaliza update.
$ws.Range("D2").Value = 44559
$ws.Range("D3").Value = 44559
$ws.Range("D4").Value = 44574
$ws.Range("D5").Value = 44574
